$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.110.70"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.814.88"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5911"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2726"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06799"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07545"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "1.820.47"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.654"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6171"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "75.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.65%  "
$ws.Range("D17").Value = "28.884.27"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.437"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.755"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.784"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06346"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.417"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.736"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.701"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.691"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.065"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.541"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6338"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.747"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01714"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.460"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "1.125.58"
$ws.Range("E40").Value = "  -8.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8729"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.36%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "1.968.55"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000114"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.583"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05501"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4540"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.282"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("E51").Value = "  -3.62%  "
